$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 6): Nombre, NIF, Email (hyperlink), Colegio
$ws.Range("A6").Value = "Victor Retortillo"
$ws.Range("B6").Value = "41310533A"
$ws.Range("C6").Value = "victor@email."
$ws.Range("D6").Value = 263

# Turn the new email cell into a mailto hyperlink, like the rows above it
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:victor@email.", "", "", "victor@email.")

# Match the visual style used by the other Email column cells (hyperlink look)
$ws.Range("C6").Style = $ws.Range("C5").Style

# Update the selected cell, as recorded in the workbook view
$null = $ws.Range("C8").Select()
